# Generate Report for Handback
# Updates the status for the "8b05e0ed-..." file (row 3) from "Ready for handoff"
# to "Handback transform failed" across all sheets (shared string), and records
# the handback/handoff file name mismatch error detail for the zh-cn and de-de
# localization sheets. Also widens the "Error Detail" column so the message is
# readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# 1) Update the shared "Status" text used by Overview!E3/F3 and the per-locale
#    sheets' Status column (C3) for the 8b05e0ed-... row.
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZh.Range("C3").Value = "Handback transform failed"
$wsDe.Range("C3").Value = "Handback transform failed"

# 2) Record the "Error Detail" for the handback/handoff filename mismatch on
#    the zh-cn and de-de sheets, row 3 (8b05e0ed-... file), column P.
$wsZh.Range("P3").Value = "Handback file name: nwygjkup.0pk is different with handoff file name: 8b05e0ed-5060-4923-9a08-aba7cbbcf29f.95af13791b706c8a36f0dec975e0fb123fcff2a8.zh-cn."
$wsDe.Range("P3").Value = "Handback file name: nwygjkup.0pk is different with handoff file name: 8b05e0ed-5060-4923-9a08-aba7cbbcf29f.95af13791b706c8a36f0dec975e0fb123fcff2a8.de-de."

# 3) Widen the "Error Detail" column (P) on both locale sheets so the new
#    message is visible (OOXML width 40, same as other wide columns).
$wsZh.Columns.Item(16).ColumnWidth = 39.17
$wsDe.Columns.Item(16).ColumnWidth = 39.17
